$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the date/time columns as text first so Excel doesn't silently
# reinterpret "08-12-2024" / "17:11:12" as date/time serial numbers.
$ws.Range("C2:C3").NumberFormat = "@"
$ws.Range("D2:D3").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 1).Value = "Dinnar Ary Nastiti"
$ws.Cells.Item(2, 2).Value = "MJM005"
$ws.Cells.Item(2, 3).Value = "08-12-2024"
$ws.Cells.Item(2, 4).Value = "17:11:12"
$ws.Cells.Item(2, 5).Value = "Alfa"

# Row 3
$ws.Cells.Item(3, 1).Value = "Dinnar Ary Nastiti"
$ws.Cells.Item(3, 2).Value = "MJM005"
$ws.Cells.Item(3, 3).Value = "08-12-2024"
$ws.Cells.Item(3, 4).Value = "17:22:06"
$ws.Cells.Item(3, 5).Value = "Alfa"

# Restore the default "Normal" style on the text-formatted cells so they
# don't carry an extra cell style compared to the original workbook.
$ws.Range("C2:D3").Style = "Normal"
